# Weekly update: insert a new "Ajo" (garlic) price record as the most
# recent entry in the Femacal de La Calera Hortaliza sheet.
#
# This pushes the previously-existing rows 196-224 down to 197-225 and
# fills the freshly opened row 196 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 196; this shifts every row
# from 196 onward down by one (old row 196 becomes 197, ..., old row
# 224 becomes 225) and copies formatting from the row being split.
$ws.Rows.Item(196).Insert()

# Populate the newly inserted row 196 with the new week's record.
$ws.Cells.Item(196, 1).Value2  = 3
$ws.Cells.Item(196, 2).Value2  = "Femacal de La Calera"
$ws.Cells.Item(196, 3).Value2  = "Coquimbo"
$ws.Cells.Item(196, 4).Value2  = 44474
$ws.Cells.Item(196, 5).Value2  = 5
$ws.Cells.Item(196, 6).Value2  = 100112003
$ws.Cells.Item(196, 7).Value2  = "Ajo"
$ws.Cells.Item(196, 8).Value2  = "Chino"
$ws.Cells.Item(196, 9).Value2  = "Primera"
$ws.Cells.Item(196, 10).Value2 = 73
$ws.Cells.Item(196, 11).Value2 = 16500
$ws.Cells.Item(196, 12).Value2 = 17000
$ws.Cells.Item(196, 13).Value2 = 16760
$ws.Cells.Item(196, 14).Value2 = "`$/caja 10 kilos"
$ws.Cells.Item(196, 15).Value2 = "China"
$ws.Cells.Item(196, 16).Value2 = 1676
$ws.Cells.Item(196, 17).Value2 = 10
$ws.Cells.Item(196, 18).Value2 = "Hortaliza"
